$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.851.43"
$ws.Range("E2").Value = "  -0.98%  "

$ws.Range("D3").Value = "2.682.84"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'599.64"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").Value = "'167.77"
$ws.Range("E6").Value = "  +4.52%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +0.83%  "

$ws.Range("D9").Value = "2.681.95"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("D10").Value = "'0.144"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("D12").Value = "'0.360"
$ws.Range("E12").Value = "  +0.04%  "

$ws.Range("D13").Value = "'5.24"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("D14").Value = "'27.98"
$ws.Range("E14").Value = "  -1.15%  "

$ws.Range("D15").Value = "3.168.82"
$ws.Range("E15").Value = "  -0.66%  "

$ws.Range("D16").Value = "'0.0000186"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").Value = "67.745.56"
$ws.Range("E17").Value = "  -1.11%  "

$ws.Range("D18").Value = "2.679.09"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").Value = "'11.78"
$ws.Range("E19").Value = "  -0.42%  "

$ws.Range("D20").Value = "'7.71"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("D21").Value = "'364.81"
$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'4.40"
$ws.Range("E22").Value = "  -2.97%  "

$ws.Range("D23").Value = "'4.85"
$ws.Range("E23").Value = "  -0.77%  "

$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = "  -3.57%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("D26").Value = "'71.03"
$ws.Range("E26").Value = "  -4.38%  "

$ws.Range("D27").Value = "'10.07"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").Value = "2.721.76"
$ws.Range("E28").Value = "  -3.91%  "

$ws.Range("E29").Value = "  -1.76%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("D31").Value = "'561.32"
$ws.Range("E31").Value = "  -4.35%  "

$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = "  -1.92%  "

$ws.Range("E33").Value = "  -3.07%  "

$ws.Range("D34").Value = "'1.94"
$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").Value = "'0.131"
$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -4.34%  "

$ws.Range("D38").Value = "'19.60"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("D39").Value = "'156.28"
$ws.Range("E39").Value = "  -2.86%  "

$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("D41").Value = "'5.34"
$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("E42").Value = "  -3.60%  "

$ws.Range("D43").Value = "'17.98"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").Value = "'2.54"
$ws.Range("E44").Value = "  -4.82%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "'40.32"
$ws.Range("E46").Value = "  -0.97%  "

$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("E47").Value = "  -5.24%  "

$ws.Range("D48").Value = "'0.593"
$ws.Range("E48").Value = "  -1.78%  "

$ws.Range("D49").Value = "'153.90"
$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("D50").Value = "'3.86"
$ws.Range("E50").Value = "  -1.75%  "

$ws.Range("E51").Value = "  -2.21%  "

# Reset number format artifacts introduced by text-forcing apostrophe prefix
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
